# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet a new (blank) column is inserted
# before the existing "Late" column, shifting the old N/O/P columns
# (Late / heading / Outstanding) one to the right. The newly inserted
# column inherits the width of its left neighbour (column M = 11
# characters) without the "best fit" flag, matching a manual
# right-click > Insert on the column header in Excel.
#
# The "Repayment schedule" tab also becomes the active/selected sheet
# (it previously was "NewLoanInput"), with cell T10 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late").
$ws.Columns("N:N").Insert() | Out-Null

# The inserted column takes on column M's width (11 characters).
$ws.Columns("N:N").ColumnWidth = 10.166666666666666

# Make "Repayment schedule" the active sheet/tab, with T10 selected
# (this also clears tabSelected on the previously active sheet).
$ws.Activate() | Out-Null
$ws.Range("T10").Select() | Out-Null
